# Add season-record columns (Wins, Losses, Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — new column headers in AC, AD, AE.
$ws.Cells.Item(1, 29).Value = "Wins"
$ws.Cells.Item(1, 30).Value = "Losses"
$ws.Cells.Item(1, 31).Value = "Ties"

# Match the bold/centered/bordered header style used by the rest of row 1
# (copy the formatting from the neighboring "Unnamed: 27" header cell).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Data rows 2-37 — season record for every player: 81 wins, 81 losses, 0 ties.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 81
    $ws.Cells.Item($r, 30).Value = 81
    $ws.Cells.Item($r, 31).Value = 0
}
